# Hyperparameter tuning for meta subclasses: update evaluation metrics
# for the kNN (row 2) and Ensemble (row 6) models.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ("kNN")
$ws.Range("B2").Value = 0.6228571428571429
$ws.Range("C2").Value = 0.6028299319727891
$ws.Range("D2").Value = 0.6228571428571429
$ws.Range("E2").Value = 0.5631802054154995
$ws.Range("J2").Value = 0.7685714285714285
$ws.Range("K2").Value = 0.8038095238095238
$ws.Range("L2").Value = 0.7685714285714285
$ws.Range("M2").Value = 0.7549641311069882
$ws.Range("N2").Value = 0.8114285714285714
$ws.Range("O2").Value = 0.8068253968253968
$ws.Range("P2").Value = 0.8114285714285714
$ws.Range("Q2").Value = 0.7898923298923299
$ws.Range("R2").Value = 0.580952380952381
$ws.Range("S2").Value = 0.5471395271395272
$ws.Range("T2").Value = 0.580952380952381
$ws.Range("U2").Value = 0.4939516192441888
$ws.Range("V2").Value = 0.580952380952381
$ws.Range("W2").Value = 0.5635031635031635
$ws.Range("X2").Value = 0.580952380952381
$ws.Range("Y2").Value = 0.4948143643422281

# Row 6 ("Ensemble")
$ws.Range("B6").Value = 0.7028571428571428
$ws.Range("C6").Value = 0.7156046176046177
$ws.Range("D6").Value = 0.7028571428571428
$ws.Range("E6").Value = 0.6605703614779245
$ws.Range("F6").Value = 0.8114285714285714
$ws.Range("G6").Value = 0.7911746031746032
$ws.Range("H6").Value = 0.8114285714285714
$ws.Range("I6").Value = 0.7760461760461761
$ws.Range("N6").Value = 0.7438095238095238
$ws.Range("O6").Value = 0.726374458874459
$ws.Range("P6").Value = 0.7438095238095238
$ws.Range("Q6").Value = 0.6975255332902393
$ws.Range("S6").Value = 0.7025649350649351
$ws.Range("U6").Value = 0.6627807225454284
$ws.Range("V6").Value = 0.7428571428571429
$ws.Range("W6").Value = 0.7273015873015873
$ws.Range("X6").Value = 0.7428571428571429
$ws.Range("Y6").Value = 0.6956553287981858
